$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (C) and P_Value (D) columns per corrected Diebold-Mariano test results

$ws.Range("C2").Value = 0.2993704532142381
$ws.Range("D2").Value = 0.7674673634114755

$ws.Range("C3").Value = 0.9669118278236888
$ws.Range("D3").Value = 0.3440980840260441

$ws.Range("C4").Value = 2.744122980905169
$ws.Range("D4").Value = 0.01184251181707374

$ws.Range("C5").Value = 3.192566322957168
$ws.Range("D5").Value = 0.004205060765871238

$ws.Range("C6").Value = 0.6639443969063273
$ws.Range("D6").Value = 0.513624615254755

$ws.Range("C7").Value = 2.789356599532246
$ws.Range("D7").Value = 0.0106906919234353

$ws.Range("C8").Value = 2.815533950627114
$ws.Range("D8").Value = 0.01007357481652327

$ws.Range("C9").Value = 2.331328304349616
$ws.Range("D9").Value = 0.0292967797309116

$ws.Range("C10").Value = 4.077472475926144
$ws.Range("D10").Value = 0.0004993021002315245

$ws.Range("C11").Value = 0.2169943424381551
$ws.Range("D11").Value = 0.830212472195039

$wb.Save()
